$d = $word.ActiveDocument

# 1. Update the CSS locator for question 2 (the 10 "tbody[class='mdc-data-table__content'] ..."
#    lines) to the new example-viewer based CSS path.
for ($i = 1; $i -le 10; $i++) {
    $old = "tbody[class='mdc-data-table__content'] tr:nth-child($i) td:nth-child(3)"
    $new = "example-viewer:nth-child(1) > div:nth-child(1) > div:nth-child(2) > table-basic-example:nth-child(1) > table:nth-child(1) > tbody:nth-child(2) > tr:nth-child($i) > td:nth-child(3)"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 2. Move the "_GoBack" bookmark from the end of the document (after the question 4 CSS
#    answer) up to the empty paragraph right after the question 2 CSS answers.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $d.Paragraphs.Item(41)
$d.Bookmarks.Add("_GoBack", $target.Range) | Out-Null
